# Update "Förändrad" (Changed) date column C for all data rows (2..210)
# from serial date 45179 (2023-09-10) to 45180 (2023-09-11).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 210; $r++) {
    $ws.Cells.Item($r, 3).Value = 45180
}
